# Insert a new data row before the existing row 280 ("Feria Lagunitas de
# Puerto Montt" / Pera sheet). This pushes the former rows 280-373 down to
# 281-374 and the sheet's used range grows from A1:T373 to A1:T374.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 280 (and everything below it) down by one row.
$ws.Rows("280:280").Insert()

# Populate the newly inserted row 280 with the new record.
$ws.Cells.Item(280, 1).Value2  = 4
$ws.Cells.Item(280, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(280, 3).Value2  = "Los Lagos"
$ws.Cells.Item(280, 4).Value2  = 44876
$ws.Cells.Item(280, 5).Value2  = 10
$ws.Cells.Item(280, 6).Value2  = "Fruta"
$ws.Cells.Item(280, 7).Value2  = 100104
$ws.Cells.Item(280, 8).Value2  = "Frutos de pepita"
$ws.Cells.Item(280, 9).Value2  = 100104005
$ws.Cells.Item(280, 10).Value2 = "Pera"
$ws.Cells.Item(280, 11).Value2 = "Packham's Triumph"
$ws.Cells.Item(280, 12).Value2 = "Primera"
$ws.Cells.Item(280, 13).Value2 = 400
$ws.Cells.Item(280, 14).Value2 = 17000
$ws.Cells.Item(280, 15).Value2 = 18000
$ws.Cells.Item(280, 16).Value2 = 17500
$ws.Cells.Item(280, 17).Value2 = "$/caja 15 kilos empedrada"
$ws.Cells.Item(280, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(280, 19).Value2 = 1167
$ws.Cells.Item(280, 20).Value2 = 15
